$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $bVal = $ws.Cells.Item($row, 2).Value2
    $cVal = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 2).Value = $cVal
    $ws.Cells.Item($row, 3).Value = $bVal
}
